# Apply the "2 dernieres corrections" update:
#  - add an extra step column (ETAPE 11 and ETAPE 12 data) to the 3 report sheets
#  - remove the now-unused last column (ETAPE 13) that was never filled in
#  - update the active sheet / selections to reflect where the author ended up

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LightHouse - Portable")
$ws2 = $wb.Worksheets.Item("LightHouse - Bureau")
$ws3 = $wb.Worksheets.Item("GTmetrix - Bureau")

# ---------------------------------------------------------------------------
# Sheet 1 : LightHouse - Portable
# ---------------------------------------------------------------------------
$ws1.Range("N4").Value = 83
$ws1.Range("O4").Value = 83

$ws1.Range("N5").Value = 96
$ws1.Range("O5").Value = 96

$ws1.Range("N6").Value = 87
$ws1.Range("O6").Value = 87

$ws1.Range("N7").Value = 89
$ws1.Range("O7").Value = 97

$ws1.Range("N4:N7").Borders.Weight = $ws1.Range("M4:M7").Borders.Weight
$ws1.Range("O4:O7").Borders.Weight = $ws1.Range("M4:M7").Borders.Weight

# drop column P entirely (it only ever held the "ETAPE 13" header, never data)
$ws1.Range("P:P").Delete()

# ---------------------------------------------------------------------------
# Sheet 2 : LightHouse - Bureau
# ---------------------------------------------------------------------------
$ws2.Range("N4").Value = 91
$ws2.Range("O4").Value = 92

$ws2.Range("N5").Value = 95
$ws2.Range("O5").Value = 95

$ws2.Range("N6").Value = 96
$ws2.Range("O6").Value = 93

$ws2.Range("N7").Value = 100
$ws2.Range("O7").Value = 100

$ws2.Range("N4:N7").Borders.Weight = $ws2.Range("M4:M7").Borders.Weight
$ws2.Range("O4:O7").Borders.Weight = $ws2.Range("M4:M7").Borders.Weight

$ws2.Range("P:P").Delete()

# ---------------------------------------------------------------------------
# Sheet 3 : GTmetrix - Bureau
# ---------------------------------------------------------------------------
$ws3.Range("N4").Value = 98
$ws3.Range("O4").Value = 99

$ws3.Range("N5").Value = 95
$ws3.Range("O5").Value = 96

$ws3.Range("N4:N5").Borders.Weight = $ws3.Range("M4:M5").Borders.Weight

$ws3.Range("P:P").Delete()

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$ws1.Range("F17").Select()
$ws2.Range("B2:O7").Select()
$ws3.Range("B2:O5").Select()

$ws3.Activate()
$ws3.Select()

$wb.Save()
